$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.152.28"
$ws.Range("E2").Value = "  +7.77%  "

$ws.Range("D3").Value = "1.585.89"
$ws.Range("E3").Value = "  +7.68%  "

$ws.Range("E4").Value = "  -0.07%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.9909"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "297.57"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +7.25%  "

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3617"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.21%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3333"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +8.34%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "41.09"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.53%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.110"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.45%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.06923"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.14%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.002"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.12%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "19.30"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.87%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.802"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +5.25%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.506"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +5.55%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.9917"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +3.02%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001059"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +3.14%  "

$ws.Range("D18").Value = "1.585.98"
$ws.Range("E18").Value = "  +7.54%  "

$ws.Range("E19").Value = "  +10.86%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "75.78"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +10.13%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "15.76"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +8.39%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.896"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +7.12%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "11.53"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +3.13%  "

$ws.Range("D24").Value = "22.159.64"
$ws.Range("E24").Value = "  +7.74%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.360"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +4.67%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.482"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +16.48%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "148.31"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +4.37%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.07"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +10.98%  "

$ws.Range("D29").Value = "1.753.32"
$ws.Range("E29").Value = "  +7.26%  "

$ws.Range("E30").Value = "  +6.38%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "3.926"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +1.00%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "5.846"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +18.27%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.9121"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +13.19%  "

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.08134"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.65%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.632"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +8.18%  "

$ws.Range("E36").Value = "  +12.06%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "5.079"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +7.44%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.231"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +1.45%  "

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05987"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +3.80%  "

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "8.245"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +11.12%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.02169"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +5.78%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.1971"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +4.96%  "

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.9911"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.91%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.5740"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +8.70%  "

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.758"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.74%  "

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "12.88"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.23%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "124.62"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +4.83%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.5539"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +6.46%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.930"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +6.19%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.06699"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +3.64%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "72.03"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +7.16%  "
